# Avance en carga financiera, monto original, monto actual
#
# The "interfaz_rel" sheet (tab 2) lists the fields of each record type
# (B01, C01, D01, ...). A new field "otros_pagado" (numeric) needs to be
# inserted right after the last B01 row ("otros"), i.e. as a brand-new
# row 39. Every row from the old row 39 onward shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a new blank row at position 39 - this pushes the old rows
# 39..50 down to 40..51, carrying their values/styles with them.
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row with the new "otros_pagado" field.
# Inserting a row copies the formatting of the row above, so column
# C/D already pick up the "label, left/valign-center" / "center/center"
# styles used by the rest of the B01 block (rows 29-38) automatically.
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = "B01"
$ws.Cells.Item(39, 3).Value = "otros_pagado"
$ws.Cells.Item(39, 4).Value = "numeric"

# Column A is a simple 1-based running id (row number - 1). Re-stamp it
# for the rows that were just shifted down so the sequence stays
# contiguous (38, 39, 40, ... 50).
for ($r = 40; $r -le 51; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Reflect the editor's final cursor position.
$ws.Activate()
[void]$ws.Range("D39").Select()
